# 🔄 MAJ automatique BRVM via GitHub Actions
# Refreshes the "Recommandations" and "Top_YTD" sheets with the latest BRVM
# market-data snapshot (titles re-ranked by "Variation Totale (%)" /
# "Progression YTD (%)", figures refreshed, two stale rows dropped).

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsTop  = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet: the list shrank from 50 to 48 titles, so drop the
# two now-unused trailing rows (50 and 51) before rewriting the remaining ones ---
$wsReco.Range("A50:A51").EntireRow.Delete() | Out-Null

# --- Recommandations sheet: refresh rows 2-49 (A:G) with the latest values ---
$wsReco.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$wsReco.Cells.Item(2, 2).Value = 0
$wsReco.Cells.Item(2, 3).Value = 8
$wsReco.Cells.Item(2, 4).Value = 3442.45
$wsReco.Cells.Item(2, 5).Value = 112.35
$wsReco.Cells.Item(2, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(2, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(3, 1).Value = "SUCRIVOIRE"
$wsReco.Cells.Item(3, 2).Value = 0
$wsReco.Cells.Item(3, 3).Value = 3
$wsReco.Cells.Item(3, 4).Value = 2760
$wsReco.Cells.Item(3, 5).Value = 945
$wsReco.Cells.Item(3, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(3, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(4, 1).Value = "SAFCA CI"
$wsReco.Cells.Item(4, 2).Value = 0
$wsReco.Cells.Item(4, 3).Value = 4
$wsReco.Cells.Item(4, 4).Value = 2755
$wsReco.Cells.Item(4, 5).Value = 690
$wsReco.Cells.Item(4, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(4, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(5, 1).Value = "CFAO MOTORS CI"
$wsReco.Cells.Item(5, 2).Value = 0
$wsReco.Cells.Item(5, 3).Value = 4
$wsReco.Cells.Item(5, 4).Value = 2695
$wsReco.Cells.Item(5, 5).Value = 665
$wsReco.Cells.Item(5, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(5, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(6, 1).Value = "BRVM - AUTRES SECTEURS"
$wsReco.Cells.Item(6, 2).Value = 0
$wsReco.Cells.Item(6, 3).Value = 4
$wsReco.Cells.Item(6, 4).Value = 2646.43
$wsReco.Cells.Item(6, 5).Value = 664.36
$wsReco.Cells.Item(6, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(6, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(7, 1).Value = "UNIWAX CI"
$wsReco.Cells.Item(7, 2).Value = 0
$wsReco.Cells.Item(7, 3).Value = 4
$wsReco.Cells.Item(7, 4).Value = 2410
$wsReco.Cells.Item(7, 5).Value = 600
$wsReco.Cells.Item(7, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(7, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(8, 1).Value = "NEI-CEDA CI"
$wsReco.Cells.Item(8, 2).Value = 0
$wsReco.Cells.Item(8, 3).Value = 4
$wsReco.Cells.Item(8, 4).Value = 2360
$wsReco.Cells.Item(8, 5).Value = 585
$wsReco.Cells.Item(8, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(8, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(9, 1).Value = "AIR LIQUIDE CI"
$wsReco.Cells.Item(9, 2).Value = 0
$wsReco.Cells.Item(9, 3).Value = 4
$wsReco.Cells.Item(9, 4).Value = 2175
$wsReco.Cells.Item(9, 5).Value = 545
$wsReco.Cells.Item(9, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(9, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(10, 1).Value = "SETAO CI"
$wsReco.Cells.Item(10, 2).Value = 0
$wsReco.Cells.Item(10, 3).Value = 4
$wsReco.Cells.Item(10, 4).Value = 2125
$wsReco.Cells.Item(10, 5).Value = 540
$wsReco.Cells.Item(10, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(10, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(11, 1).Value = "BRVM - TRANSPORT"
$wsReco.Cells.Item(11, 2).Value = 0
$wsReco.Cells.Item(11, 3).Value = 4
$wsReco.Cells.Item(11, 4).Value = 1597.12
$wsReco.Cells.Item(11, 5).Value = 397.79
$wsReco.Cells.Item(11, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(11, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(12, 1).Value = "BRVM - DISTRIBUTION"
$wsReco.Cells.Item(12, 2).Value = 0
$wsReco.Cells.Item(12, 3).Value = 4
$wsReco.Cells.Item(12, 4).Value = 1496.11
$wsReco.Cells.Item(12, 5).Value = 369.92
$wsReco.Cells.Item(12, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(12, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(13, 1).Value = "BRVM - AGRICULTURE"
$wsReco.Cells.Item(13, 2).Value = 0
$wsReco.Cells.Item(13, 3).Value = 4
$wsReco.Cells.Item(13, 4).Value = 1314.66
$wsReco.Cells.Item(13, 5).Value = 335.77
$wsReco.Cells.Item(13, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(13, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(14, 1).Value = "BRVM - INDUSTRIE"
$wsReco.Cells.Item(14, 2).Value = 0
$wsReco.Cells.Item(14, 3).Value = 4
$wsReco.Cells.Item(14, 4).Value = 775.08
$wsReco.Cells.Item(14, 5).Value = 191.98
$wsReco.Cells.Item(14, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(14, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(15, 1).Value = "BRVM-PRINCIPAL"
$wsReco.Cells.Item(15, 2).Value = 0
$wsReco.Cells.Item(15, 3).Value = 4
$wsReco.Cells.Item(15, 4).Value = 714.02
$wsReco.Cells.Item(15, 5).Value = 179.02
$wsReco.Cells.Item(15, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(15, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(16, 1).Value = "BRVM - CONSOMMATION DE BASE"
$wsReco.Cells.Item(16, 2).Value = 0
$wsReco.Cells.Item(16, 3).Value = 4
$wsReco.Cells.Item(16, 4).Value = 675.31
$wsReco.Cells.Item(16, 5).Value = 169.05
$wsReco.Cells.Item(16, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(16, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(17, 1).Value = "BRVM - INDUSTRIELS"
$wsReco.Cells.Item(17, 2).Value = 0
$wsReco.Cells.Item(17, 3).Value = 4
$wsReco.Cells.Item(17, 4).Value = 607.29
$wsReco.Cells.Item(17, 5).Value = 149.5
$wsReco.Cells.Item(17, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(17, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(18, 1).Value = "BRVM-PRESTIGE"
$wsReco.Cells.Item(18, 2).Value = 0
$wsReco.Cells.Item(18, 3).Value = 4
$wsReco.Cells.Item(18, 4).Value = 520.82
$wsReco.Cells.Item(18, 5).Value = 131.31
$wsReco.Cells.Item(18, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(18, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(19, 1).Value = "BRVM - FINANCES"
$wsReco.Cells.Item(19, 2).Value = 0
$wsReco.Cells.Item(19, 3).Value = 4
$wsReco.Cells.Item(19, 4).Value = 491.27
$wsReco.Cells.Item(19, 5).Value = 124.56
$wsReco.Cells.Item(19, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(19, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(20, 1).Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Cells.Item(20, 2).Value = 0
$wsReco.Cells.Item(20, 3).Value = 4
$wsReco.Cells.Item(20, 4).Value = 482.81
$wsReco.Cells.Item(20, 5).Value = 122.41
$wsReco.Cells.Item(20, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(20, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(21, 1).Value = "BRVM - ENERGIE"
$wsReco.Cells.Item(21, 2).Value = 0
$wsReco.Cells.Item(21, 3).Value = 4
$wsReco.Cells.Item(21, 4).Value = 448.17
$wsReco.Cells.Item(21, 5).Value = 111.34
$wsReco.Cells.Item(21, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(21, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(22, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Cells.Item(22, 2).Value = 0
$wsReco.Cells.Item(22, 3).Value = 4
$wsReco.Cells.Item(22, 4).Value = 429.2
$wsReco.Cells.Item(22, 5).Value = 106.61
$wsReco.Cells.Item(22, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(22, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(23, 1).Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Cells.Item(23, 2).Value = 0
$wsReco.Cells.Item(23, 3).Value = 4
$wsReco.Cells.Item(23, 4).Value = 390.24
$wsReco.Cells.Item(23, 5).Value = 97.23
$wsReco.Cells.Item(23, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(23, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(24, 1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$wsReco.Cells.Item(24, 2).Value = 2
$wsReco.Cells.Item(24, 3).Value = 0
$wsReco.Cells.Item(24, 4).Value = 8.91
$wsReco.Cells.Item(24, 5).Value = 4.63
$wsReco.Cells.Item(24, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(24, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(25, 1).Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Cells.Item(25, 2).Value = 1
$wsReco.Cells.Item(25, 3).Value = 0
$wsReco.Cells.Item(25, 4).Value = 7.5
$wsReco.Cells.Item(25, 5).Value = 7.5
$wsReco.Cells.Item(25, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(25, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(26, 1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$wsReco.Cells.Item(26, 2).Value = 2
$wsReco.Cells.Item(26, 3).Value = 0
$wsReco.Cells.Item(26, 4).Value = 7.44
$wsReco.Cells.Item(26, 5).Value = 4.3
$wsReco.Cells.Item(26, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(26, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(27, 1).Value = "SOLIBRA CI (SLBC)"
$wsReco.Cells.Item(27, 2).Value = 1
$wsReco.Cells.Item(27, 3).Value = 0
$wsReco.Cells.Item(27, 4).Value = 7.44
$wsReco.Cells.Item(27, 5).Value = 7.44
$wsReco.Cells.Item(27, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(27, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(28, 1).Value = "AIR LIQUIDE CI (SIVC)"
$wsReco.Cells.Item(28, 2).Value = 1
$wsReco.Cells.Item(28, 3).Value = 0
$wsReco.Cells.Item(28, 4).Value = 6.93
$wsReco.Cells.Item(28, 5).Value = 6.93
$wsReco.Cells.Item(28, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(28, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(29, 1).Value = "SOGB CI (SOGC)"
$wsReco.Cells.Item(29, 2).Value = 1
$wsReco.Cells.Item(29, 3).Value = 0
$wsReco.Cells.Item(29, 4).Value = 6.45
$wsReco.Cells.Item(29, 5).Value = 6.45
$wsReco.Cells.Item(29, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(29, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(30, 1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$wsReco.Cells.Item(30, 2).Value = 2
$wsReco.Cells.Item(30, 3).Value = 1
$wsReco.Cells.Item(30, 4).Value = 5.41
$wsReco.Cells.Item(30, 5).Value = -7.47
$wsReco.Cells.Item(30, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(30, 7).Value = "👀 À surveiller"
$wsReco.Cells.Item(31, 1).Value = "BERNABE CI (BNBC)"
$wsReco.Cells.Item(31, 2).Value = 2
$wsReco.Cells.Item(31, 3).Value = 1
$wsReco.Cells.Item(31, 4).Value = 5.35
$wsReco.Cells.Item(31, 5).Value = 3.02
$wsReco.Cells.Item(31, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(31, 7).Value = "👀 À surveiller"
$wsReco.Cells.Item(32, 1).Value = "SMB CI (SMBC)"
$wsReco.Cells.Item(32, 2).Value = 2
$wsReco.Cells.Item(32, 3).Value = 0
$wsReco.Cells.Item(32, 4).Value = 3.74
$wsReco.Cells.Item(32, 5).Value = 3.74
$wsReco.Cells.Item(32, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(32, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(33, 1).Value = "UNILEVER CI (UNLC)"
$wsReco.Cells.Item(33, 2).Value = 1
$wsReco.Cells.Item(33, 3).Value = 0
$wsReco.Cells.Item(33, 4).Value = 3.59
$wsReco.Cells.Item(33, 5).Value = 3.59
$wsReco.Cells.Item(33, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(33, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(34, 1).Value = "SUCRIVOIRE (SCRC)"
$wsReco.Cells.Item(34, 2).Value = 1
$wsReco.Cells.Item(34, 3).Value = 0
$wsReco.Cells.Item(34, 4).Value = 3.09
$wsReco.Cells.Item(34, 5).Value = 3.09
$wsReco.Cells.Item(34, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(34, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(35, 1).Value = "SETAO CI (STAC)"
$wsReco.Cells.Item(35, 2).Value = 1
$wsReco.Cells.Item(35, 3).Value = 1
$wsReco.Cells.Item(35, 4).Value = 3.02
$wsReco.Cells.Item(35, 5).Value = 5.88
$wsReco.Cells.Item(35, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(35, 7).Value = "👀 À surveiller"
$wsReco.Cells.Item(36, 1).Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Cells.Item(36, 2).Value = 1
$wsReco.Cells.Item(36, 3).Value = 1
$wsReco.Cells.Item(36, 4).Value = 2.41
$wsReco.Cells.Item(36, 5).Value = -2.21
$wsReco.Cells.Item(36, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(36, 7).Value = "👀 À surveiller"
$wsReco.Cells.Item(37, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Cells.Item(37, 2).Value = 2
$wsReco.Cells.Item(37, 3).Value = 1
$wsReco.Cells.Item(37, 4).Value = 0.32
$wsReco.Cells.Item(37, 5).Value = 5.88
$wsReco.Cells.Item(37, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(37, 7).Value = "👀 À surveiller"
$wsReco.Cells.Item(38, 1).Value = "SAFCA CI (SAFC)"
$wsReco.Cells.Item(38, 2).Value = 1
$wsReco.Cells.Item(38, 3).Value = 1
$wsReco.Cells.Item(38, 4).Value = 0.2
$wsReco.Cells.Item(38, 5).Value = 4.55
$wsReco.Cells.Item(38, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(38, 7).Value = "👀 À surveiller"
$wsReco.Cells.Item(39, 1).Value = "TOTAL"
$wsReco.Cells.Item(39, 2).Value = 0
$wsReco.Cells.Item(39, 3).Value = 4
$wsReco.Cells.Item(39, 4).Value = 0
$wsReco.Cells.Item(39, 5).Value = 0
$wsReco.Cells.Item(39, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(39, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(40, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Cells.Item(40, 2).Value = 1
$wsReco.Cells.Item(40, 3).Value = 1
$wsReco.Cells.Item(40, 4).Value = -0.17
$wsReco.Cells.Item(40, 5).Value = 3.44
$wsReco.Cells.Item(40, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(40, 7).Value = "👀 À surveiller"
$wsReco.Cells.Item(41, 1).Value = "BICI CI (BICC)"
$wsReco.Cells.Item(41, 2).Value = 0
$wsReco.Cells.Item(41, 3).Value = 1
$wsReco.Cells.Item(41, 4).Value = -1.78
$wsReco.Cells.Item(41, 5).Value = -1.78
$wsReco.Cells.Item(41, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(41, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(42, 1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(42, 2).Value = 0
$wsReco.Cells.Item(42, 3).Value = 1
$wsReco.Cells.Item(42, 4).Value = -2.42
$wsReco.Cells.Item(42, 5).Value = -2.42
$wsReco.Cells.Item(42, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(42, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(43, 1).Value = "NEI-CEDA CI (NEIC)"
$wsReco.Cells.Item(43, 2).Value = 0
$wsReco.Cells.Item(43, 3).Value = 1
$wsReco.Cells.Item(43, 4).Value = -2.5
$wsReco.Cells.Item(43, 5).Value = -2.5
$wsReco.Cells.Item(43, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(43, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(44, 1).Value = "BANK OF AFRICA CI (BOAC)"
$wsReco.Cells.Item(44, 2).Value = 0
$wsReco.Cells.Item(44, 3).Value = 1
$wsReco.Cells.Item(44, 4).Value = -2.78
$wsReco.Cells.Item(44, 5).Value = -2.78
$wsReco.Cells.Item(44, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(44, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(45, 1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$wsReco.Cells.Item(45, 2).Value = 0
$wsReco.Cells.Item(45, 3).Value = 2
$wsReco.Cells.Item(45, 4).Value = -5.46
$wsReco.Cells.Item(45, 5).Value = -3.64
$wsReco.Cells.Item(45, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(45, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(46, 1).Value = "FILTISAC CI (FTSC)"
$wsReco.Cells.Item(46, 2).Value = 0
$wsReco.Cells.Item(46, 3).Value = 2
$wsReco.Cells.Item(46, 4).Value = -5.83
$wsReco.Cells.Item(46, 5).Value = -3.55
$wsReco.Cells.Item(46, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(46, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(47, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$wsReco.Cells.Item(47, 2).Value = 0
$wsReco.Cells.Item(47, 3).Value = 1
$wsReco.Cells.Item(47, 4).Value = -7.21
$wsReco.Cells.Item(47, 5).Value = -7.21
$wsReco.Cells.Item(47, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(47, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(48, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$wsReco.Cells.Item(48, 2).Value = 0
$wsReco.Cells.Item(48, 3).Value = 2
$wsReco.Cells.Item(48, 4).Value = -7.64
$wsReco.Cells.Item(48, 5).Value = -4.18
$wsReco.Cells.Item(48, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(48, 7).Value = "➖ Neutre"
$wsReco.Cells.Item(49, 1).Value = "ONATEL BF (ONTBF)"
$wsReco.Cells.Item(49, 2).Value = 0
$wsReco.Cells.Item(49, 3).Value = 2
$wsReco.Cells.Item(49, 4).Value = -12.53
$wsReco.Cells.Item(49, 5).Value = -7.5
$wsReco.Cells.Item(49, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(49, 7).Value = "➖ Neutre"

# --- Top_YTD sheet: refresh rows 2-11 (A:B) with the latest values ---
$wsTop.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$wsTop.Cells.Item(2, 2).Value = 10439245.11
$wsTop.Cells.Item(3, 1).Value = "SAFCA CI"
$wsTop.Cells.Item(3, 2).Value = 386935.62
$wsTop.Cells.Item(4, 1).Value = "CFAO MOTORS CI"
$wsTop.Cells.Item(4, 2).Value = 358278.02
$wsTop.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsTop.Cells.Item(5, 2).Value = 336336.48
$wsTop.Cells.Item(6, 1).Value = "UNIWAX CI"
$wsTop.Cells.Item(6, 2).Value = 243332
$wsTop.Cells.Item(7, 1).Value = "NEI-CEDA CI"
$wsTop.Cells.Item(7, 2).Value = 226511.7
$wsTop.Cells.Item(8, 1).Value = "AIR LIQUIDE CI"
$wsTop.Cells.Item(8, 2).Value = 171635.12
$wsTop.Cells.Item(9, 1).Value = "SETAO CI"
$wsTop.Cells.Item(9, 2).Value = 158539.04
$wsTop.Cells.Item(10, 1).Value = "SUCRIVOIRE"
$wsTop.Cells.Item(10, 2).Value = 105894.35
$wsTop.Cells.Item(11, 1).Value = "BRVM - TRANSPORT"
$wsTop.Cells.Item(11, 2).Value = 62024.61

